$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (the extra LR-pair rows no longer present)
$ws.Rows("3:4").Delete()

# Update row 2 values per new TPM data
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.213927
$ws.Range("N2").Value = 0.641781
$ws.Range("Q2").Value = 0.330596653226
$ws.Range("R2").Value = 2.975369879034
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
